# Updates the cryptocurrency table's Price (D) and Volume(1h) (E) columns
# on the active worksheet, matching the latest scrape of coinranking.com data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as text (not a number),
# matching the workbook's existing inline-string cell type, without leaving
# behind any stray number-format/style changes on the cell.
function Set-TextValue {
    param($Worksheet, $CellRef, $TextValue)
    $Worksheet.Range($CellRef).NumberFormat = "@"
    $Worksheet.Range($CellRef).Value = $TextValue
    $Worksheet.Range($CellRef).ClearFormats()
}

# Plain text updates (percentages and prices that Excel won't mistake for numbers)
$plainUpdates = @{
    "D2" = '27.515.56'
    "E2" = '  -1.17%  '
    "D3" = '1.591.28'
    "E3" = '  -1.80%  '
    "E4" = '  +0.25%  '
    "E5" = '  -1.49%  '
    "E6" = '  -3.99%  '
    "E7" = '  +0.29%  '
    "E8" = '  -4.60%  '
    "E9" = '  -1.94%  '
    "E10" = '  -3.55%  '
    "E11" = '  -0.57%  '
    "D12" = '1.817.94'
    "E12" = '  -1.70%  '
    "D13" = '1.616.92'
    "E13" = '  -0.17%  '
    "E14" = '  -4.17%  '
    "E15" = '  -4.04%  '
    "E16" = '  -3.04%  '
    "D17" = '27.498.20'
    "E17" = '  -1.20%  '
    "E18" = '  -4.89%  '
    "E19" = '  -3.23%  '
    "D20" = "0.0$([char]0x2083)0690"
    "E20" = '  -4.03%  '
    "E22" = '  -2.64%  '
    "E23" = '  -3.92%  '
    "E24" = '  -1.31%  '
    "E25" = '  +0.43%  '
    "E26" = '  +0.25%  '
    "E27" = '  -2.55%  '
    "E28" = '  -2.90%  '
    "E29" = '  -4.58%  '
    "E30" = '  -1.34%  '
    "E31" = '  -2.82%  '
    "E32" = '  -4.05%  '
    "D33" = '1.354.21'
    "E33" = '  -1.97%  '
    "E34" = '  -4.11%  '
    "E35" = '  -2.21%  '
    "E36" = '  -0.81%  '
    "E37" = '  -4.12%  '
    "E38" = '  -2.79%  '
    "E39" = '  -2.98%  '
    "E40" = '  -3.83%  '
    "E41" = '  +0.28%  '
    "E42" = '  -3.84%  '
    "E43" = '  -2.29%  '
    "E44" = '  -2.66%  '
    "E45" = '  -4.55%  '
    "D46" = '1.728.74'
    "E46" = '  -1.72%  '
    "E47" = '  -3.05%  '
    "E48" = '  +0.16%  '
    "D49" = "0.0$([char]0x2087)0981"
    "E49" = '  -5.95%  '
    "E50" = '  -4.46%  '
    "E51" = '  -1.26%  '
}
foreach ($cellRef in $plainUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainUpdates[$cellRef]
}

# Price updates that look numeric and must be forced to stay as text
Set-TextValue $ws "D5" '207.17'
Set-TextValue $ws "D6" '0.500'
Set-TextValue $ws "D8" '22.20'
Set-TextValue $ws "D9" '0.251'
Set-TextValue $ws "D11" '0.0873'
Set-TextValue $ws "D15" '0.536'
Set-TextValue $ws "D16" '63.17'
Set-TextValue $ws "D18" '217.27'
Set-TextValue $ws "D19" '7.35'
Set-TextValue $ws "D25" '154.56'
Set-TextValue $ws "D27" '6.70'
Set-TextValue $ws "D28" '15.00'
Set-TextValue $ws "D32" '3.28'
Set-TextValue $ws "D37" '0.954'
Set-TextValue $ws "D40" '0.811'
Set-TextValue $ws "D43" '5.35'
Set-TextValue $ws "D44" '63.85'
Set-TextValue $ws "D45" '1.75'
Set-TextValue $ws "D47" '2.09'
Set-TextValue $ws "D48" '87.69'
Set-TextValue $ws "D51" '0.0496'

